$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new unit-test sample dates (same style/format as the existing date cells)
# in columns B and C for rows 2,4,5,6,7,8,9 and 14 (value = 8/7/2024, serial 45511)
$dateRows = @(2, 4, 5, 6, 7, 8, 9, 14)
foreach ($r in $dateRows) {
    $ws.Range("B3").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("B$r").Value = 45511

    $ws.Range("C3").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("C$r").Value = 45511
}

# Add conversion-factor formula for Temperature Gradient (row 27), matching the
# pattern used by the neighboring "Factors" column cells (P17:P26)
$ws.Range("P27").Formula = "=273.15*9/5"

# Update the active selection to match the saved view (scroll to row 25, select G18)
$ws.Range("A25").Select() | Out-Null
$ws.Range("G18").Select() | Out-Null
